$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: the trailing empty paragraph that currently carries bold
# paragraph-mark formatting (originally right before the two trailing
# blank paragraphs) loses that formatting and becomes a plain empty
# paragraph. Insert a fresh (non-bold) paragraph by splitting off the
# following plain paragraph, then delete the old bold-marked one.
# ------------------------------------------------------------------
$plainAfter = $d.Paragraphs.Item(15)
$plainAfter.Range.InsertParagraphBefore()
$oldBoldEmpty = $d.Paragraphs.Item(14)
$oldBoldEmpty.Range.Delete()

# ------------------------------------------------------------------
# Step 2: before ". RNF004:" insert a blank paragraph followed by a
# bold "Autenticação" heading paragraph.
# ------------------------------------------------------------------
$rnf004 = $d.Paragraphs.Item(10)
$rnf004.Range.InsertParagraphBefore()
$rnf004 = $d.Paragraphs.Item(11)
$rnf004.Range.InsertParagraphBefore()
$autenticacao = $d.Paragraphs.Item(11)
$autenticacao.Range.Text = "Autenticação"
$autenticacao.Range.Font.Bold = 1
$autenticacao.Range.Font.BoldBi = 1

# ------------------------------------------------------------------
# Step 3: before ". RNF001:" insert a bold "Estrutura" heading
# paragraph.
# ------------------------------------------------------------------
$rnf001 = $d.Paragraphs.Item(4)
$rnf001.Range.InsertParagraphBefore()
$estrutura = $d.Paragraphs.Item(4)
$estrutura.Range.Text = "Estrutura"
$estrutura.Range.Font.Bold = 1
$estrutura.Range.Font.BoldBi = 1

# ------------------------------------------------------------------
# Step 4: remove the now-redundant blank paragraph that used to sit
# between the two leading blank paragraphs and ". RNF001:".
# ------------------------------------------------------------------
$extraBlank = $d.Paragraphs.Item(3)
$extraBlank.Range.Delete()

Write-Output "done"
